# Update "想去人数" (F column) figures on the "展览" and "全部类型" sheets
# to reflect the latest generated output (commit 456a3b4).

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F5").Value = 5148
$wsExhibit.Range("F7").Value = 10029
$wsExhibit.Range("F9").Value = 554
$wsExhibit.Range("F10").Value = 100
$wsExhibit.Range("F11").Value = 51
$wsExhibit.Range("F12").Value = 761

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F7").Value = 5148
$wsAll.Range("F10").Value = 10029
$wsAll.Range("F12").Value = 554
$wsAll.Range("F13").Value = 100
$wsAll.Range("F16").Value = 51
$wsAll.Range("F17").Value = 761
